$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(34, 8).Value = 5000
$ws.Cells.Item(34, 9).Value = 5000
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 5000
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = -4797

$ws.Cells.Item(36, 8).Value = 5000
$ws.Cells.Item(36, 9).Value = 5000
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 5000
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -4285

$ws.Cells.Item(94, 8).Value = 6775.75
$ws.Cells.Item(94, 9).Value = 7300.273
$ws.Cells.Item(94, 10).Value = 1006
$ws.Cells.Item(94, 11).Value = 7300.273
$ws.Cells.Item(94, 12).Value = 1006
$ws.Cells.Item(94, 13).Value = -6849.273

$ws.Cells.Item(106, 8).Value = 1896.6666
$ws.Cells.Item(106, 9).Value = 1331.7307
$ws.Cells.Item(106, 10).Value = 3995
$ws.Cells.Item(106, 11).Value = 1331.7307
$ws.Cells.Item(106, 12).Value = 3995
$ws.Cells.Item(106, 13).Value = -700.7307000000001

$ws.Cells.Item(111, 8).Value = 1025.2142
$ws.Cells.Item(111, 9).Value = 813.6667
$ws.Cells.Item(111, 10).Value = 1183.875
$ws.Cells.Item(111, 11).Value = 2441.0001
$ws.Cells.Item(111, 12).Value = 3551.625
$ws.Cells.Item(111, 13).Value = 625.9998999999998

$ws.Cells.Item(132, 8).Value = 10080.228
$ws.Cells.Item(132, 9).Value = 10379.645
$ws.Cells.Item(132, 10).Value = 2495
$ws.Cells.Item(132, 11).Value = 31138.935
$ws.Cells.Item(132, 12).Value = 7485
$ws.Cells.Item(132, 13).Value = -28608.935

$ws.Cells.Item(137, 8).Value = 7878.415
$ws.Cells.Item(137, 9).Value = 13579.88
$ws.Cells.Item(137, 10).Value = 2787.8215
$ws.Cells.Item(137, 11).Value = 40739.64
$ws.Cells.Item(137, 12).Value = 8363.4645
$ws.Cells.Item(137, 13).Value = -38189.64

$ws.Cells.Item(138, 8).Value = 2748.7188
$ws.Cells.Item(138, 9).Value = 1689.6364
$ws.Cells.Item(138, 10).Value = 5078.7
$ws.Cells.Item(138, 11).Value = 5068.9092
$ws.Cells.Item(138, 12).Value = 15236.1
$ws.Cells.Item(138, 13).Value = 71.09079999999994
$ws.Cells.Item(138, 14).Value = -25516.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14679.488
$ws.Cells.Item(32, 9).Value = 14905.321
$ws.Cells.Item(32, 10).Value = 5194.5
$ws.Cells.Item(32, 11).Value = 14905.321
$ws.Cells.Item(32, 12).Value = 5194.5
$ws.Cells.Item(32, 13).Value = -14618.321
$ws.Cells.Item(32, 14).Value = -5768.5

$ws.Cells.Item(45, 8).Value = 3575.2144
$ws.Cells.Item(45, 9).Value = 2217.8
$ws.Cells.Item(45, 10).Value = 5141.4614
$ws.Cells.Item(45, 11).Value = 2217.8
$ws.Cells.Item(45, 12).Value = 5141.4614
$ws.Cells.Item(45, 13).Value = -1840.8

$ws.Cells.Item(74, 8).Value = 262083.39
$ws.Cells.Item(74, 9).Value = 286903.53
$ws.Cells.Item(74, 10).Value = 1472
$ws.Cells.Item(74, 11).Value = 286903.53
$ws.Cells.Item(74, 12).Value = 1472
$ws.Cells.Item(74, 13).Value = -286029.53
$ws.Cells.Item(74, 14).Value = -3220

$ws.Cells.Item(77, 8).Value = 262083.39
$ws.Cells.Item(77, 9).Value = 286903.53
$ws.Cells.Item(77, 10).Value = 1472
$ws.Cells.Item(77, 11).Value = 1434517.65
$ws.Cells.Item(77, 12).Value = 7360
$ws.Cells.Item(77, 13).Value = -1430149.65
$ws.Cells.Item(77, 14).Value = -16096

$ws.Cells.Item(97, 8).Value = 2235.8057
$ws.Cells.Item(97, 9).Value = 1639.5264
$ws.Cells.Item(97, 10).Value = 2902.2354
$ws.Cells.Item(97, 11).Value = 1639.5264
$ws.Cells.Item(97, 12).Value = 2902.2354
$ws.Cells.Item(97, 13).Value = -1143.5264
$ws.Cells.Item(97, 14).Value = -3894.2354

$ws.Cells.Item(108, 8).Value = 50000
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 50000
$ws.Cells.Item(108, 11).Value = 0
$ws.Cells.Item(108, 12).Value = 50000
$ws.Cells.Item(108, 14).Value = -57680

$ws.Cells.Item(122, 8).Value = 2354.5925
$ws.Cells.Item(122, 9).Value = 2221.5908
$ws.Cells.Item(122, 10).Value = 2939.8
$ws.Cells.Item(122, 11).Value = 6664.7724
$ws.Cells.Item(122, 12).Value = 8819.400000000001
$ws.Cells.Item(122, 13).Value = -4214.7724

$ws.Cells.Item(132, 8).Value = 1801.3695
$ws.Cells.Item(132, 9).Value = 1103.6897
$ws.Cells.Item(132, 10).Value = 2991.5293
$ws.Cells.Item(132, 11).Value = 3311.0691
$ws.Cells.Item(132, 12).Value = 8974.5879
$ws.Cells.Item(132, 13).Value = -781.0690999999997
$ws.Cells.Item(132, 14).Value = -14034.5879

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(40, 8).Value = 50000
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 50000
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 50000
$ws.Cells.Item(40, 14).Value = -50530

$ws.Cells.Item(94, 8).Value = 1289.2941
$ws.Cells.Item(94, 9).Value = 1348.1666
$ws.Cells.Item(94, 10).Value = 1148
$ws.Cells.Item(94, 11).Value = 1348.1666
$ws.Cells.Item(94, 12).Value = 1148
$ws.Cells.Item(94, 13).Value = -897.1666
$ws.Cells.Item(94, 14).Value = -2050

$ws.Cells.Item(96, 8).Value = 10806.667
$ws.Cells.Item(96, 9).Value = 15000
$ws.Cells.Item(96, 10).Value = 8710
$ws.Cells.Item(96, 11).Value = 15000
$ws.Cells.Item(96, 12).Value = 8710
$ws.Cells.Item(96, 13).Value = -12254
$ws.Cells.Item(96, 14).Value = -14202

$ws.Cells.Item(134, 8).Value = 1754.738
$ws.Cells.Item(134, 9).Value = 1773.1464
$ws.Cells.Item(134, 10).Value = 1000
$ws.Cells.Item(134, 11).Value = 5319.439200000001
$ws.Cells.Item(134, 12).Value = 3000
$ws.Cells.Item(134, 13).Value = -2784.439200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1962843.5
$ws.Cells.Item(31, 9).Value = 2633527.8
$ws.Cells.Item(31, 10).Value = 2381.923
$ws.Cells.Item(31, 11).Value = 2633527.8
$ws.Cells.Item(31, 12).Value = 2381.923
$ws.Cells.Item(31, 13).Value = -2633232.8
$ws.Cells.Item(31, 14).Value = -2971.923

$ws.Cells.Item(34, 8).Value = 1962843.5
$ws.Cells.Item(34, 9).Value = 2633527.8
$ws.Cells.Item(34, 10).Value = 2381.923
$ws.Cells.Item(34, 11).Value = 2633527.8
$ws.Cells.Item(34, 12).Value = 2381.923
$ws.Cells.Item(34, 13).Value = -2633325.8
$ws.Cells.Item(34, 14).Value = -2785.923

$ws.Cells.Item(62, 8).Value = 4980.24
$ws.Cells.Item(62, 9).Value = 2046.2307
$ws.Cells.Item(62, 10).Value = 8158.75
$ws.Cells.Item(62, 11).Value = 2046.2307
$ws.Cells.Item(62, 12).Value = 8158.75
$ws.Cells.Item(62, 13).Value = -1422.2307
$ws.Cells.Item(62, 14).Value = -9406.75

$ws.Cells.Item(65, 8).Value = 4980.24
$ws.Cells.Item(65, 9).Value = 2046.2307
$ws.Cells.Item(65, 10).Value = 8158.75
$ws.Cells.Item(65, 11).Value = 10231.1535
$ws.Cells.Item(65, 12).Value = 40793.75
$ws.Cells.Item(65, 13).Value = -7111.1535
$ws.Cells.Item(65, 14).Value = -47033.75

$ws.Cells.Item(94, 8).Value = 1141.1111
$ws.Cells.Item(94, 9).Value = 1596.3334
$ws.Cells.Item(94, 10).Value = 913.5
$ws.Cells.Item(94, 11).Value = 1596.3334
$ws.Cells.Item(94, 12).Value = 913.5
$ws.Cells.Item(94, 13).Value = -1145.3334

$ws.Cells.Item(109, 8).Value = 39994
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(109, 10).Value = 39994
$ws.Cells.Item(109, 11).Value = 0
$ws.Cells.Item(109, 12).Value = 39994
$ws.Cells.Item(109, 14).Value = -42074

$ws.Cells.Item(122, 8).Value = 1446.9
$ws.Cells.Item(122, 9).Value = 1446.9
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 4340.700000000001
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -1890.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 307
$ws.Cells.Item(12, 9).Value = 151
$ws.Cells.Item(12, 10).Value = 369.4
$ws.Cells.Item(12, 11).Value = 453
$ws.Cells.Item(12, 12).Value = 1108.2
$ws.Cells.Item(12, 13).Value = -280
$ws.Cells.Item(12, 14).Value = -1454.2

$ws.Cells.Item(60, 8).Value = 4023.6155
$ws.Cells.Item(60, 9).Value = 562.8
$ws.Cells.Item(60, 10).Value = 4847.619
$ws.Cells.Item(60, 11).Value = 1688.4
$ws.Cells.Item(60, 12).Value = 14542.857
$ws.Cells.Item(60, 13).Value = -1437.4
$ws.Cells.Item(60, 14).Value = -15044.857

$ws.Cells.Item(80, 8).Value = 4991.6665
$ws.Cells.Item(80, 9).Value = 4900
$ws.Cells.Item(80, 10).Value = 5000
$ws.Cells.Item(80, 11).Value = 14700
$ws.Cells.Item(80, 12).Value = 15000
$ws.Cells.Item(80, 13).Value = -13764

$ws.Cells.Item(83, 8).Value = 4991.6665
$ws.Cells.Item(83, 9).Value = 4900
$ws.Cells.Item(83, 10).Value = 5000
$ws.Cells.Item(83, 11).Value = 44100
$ws.Cells.Item(83, 12).Value = 45000
$ws.Cells.Item(83, 13).Value = -39420

$ws.Cells.Item(113, 8).Value = 2264.1428
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 2264.1428
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 6792.428400000001
$ws.Cells.Item(113, 14).Value = -11132.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 3666.6667
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 3666.6667
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 3666.6667
$ws.Cells.Item(22, 14).Value = -4724.6667

$ws.Cells.Item(122, 8).Value = 1883.64
$ws.Cells.Item(122, 9).Value = 2011.7778
$ws.Cells.Item(122, 10).Value = 1554.1428
$ws.Cells.Item(122, 11).Value = 6035.3334
$ws.Cells.Item(122, 12).Value = 4662.428400000001
$ws.Cells.Item(122, 13).Value = -3585.3334
$ws.Cells.Item(122, 14).Value = -9562.428400000001

$ws.Cells.Item(126, 8).Value = 2625
$ws.Cells.Item(126, 9).Value = 2670.9285
$ws.Cells.Item(126, 10).Value = 2553.5557
$ws.Cells.Item(126, 11).Value = 8012.7855
$ws.Cells.Item(126, 12).Value = 7660.6671
$ws.Cells.Item(126, 13).Value = -5542.7855
$ws.Cells.Item(126, 14).Value = -12600.6671

$ws.Cells.Item(132, 8).Value = 5000
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 14).Value = -20060
$ws.Cells.Item(132, 13).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 7319.689
$ws.Cells.Item(82, 9).Value = 5975.2
$ws.Cells.Item(82, 10).Value = 10008.667
$ws.Cells.Item(82, 11).Value = 5975.2
$ws.Cells.Item(82, 12).Value = 10008.667
$ws.Cells.Item(82, 13).Value = -5614.2
$ws.Cells.Item(82, 14).Value = -10730.667

$ws.Cells.Item(85, 8).Value = 7319.689
$ws.Cells.Item(85, 9).Value = 5975.2
$ws.Cells.Item(85, 10).Value = 10008.667
$ws.Cells.Item(85, 11).Value = 5975.2
$ws.Cells.Item(85, 12).Value = 10008.667
$ws.Cells.Item(85, 13).Value = -4727.2
$ws.Cells.Item(85, 14).Value = -12504.667

$ws.Cells.Item(93, 8).Value = 558549.2
$ws.Cells.Item(93, 9).Value = 771732.0600000001
$ws.Cells.Item(93, 10).Value = 4273.6
$ws.Cells.Item(93, 11).Value = 771732.0600000001
$ws.Cells.Item(93, 12).Value = 4273.6
$ws.Cells.Item(93, 13).Value = -770484.0600000001

$ws.Cells.Item(122, 8).Value = 9150.549000000001
$ws.Cells.Item(122, 9).Value = 10182.842
$ws.Cells.Item(122, 10).Value = 7516.0835
$ws.Cells.Item(122, 11).Value = 30548.526
$ws.Cells.Item(122, 12).Value = 22548.2505
$ws.Cells.Item(122, 13).Value = -28098.526

$ws.Cells.Item(132, 8).Value = 5160.2856
$ws.Cells.Item(132, 9).Value = 5049.5
$ws.Cells.Item(132, 10).Value = 5825
$ws.Cells.Item(132, 11).Value = 15148.5
$ws.Cells.Item(132, 12).Value = 17475
$ws.Cells.Item(132, 13).Value = -12618.5
$ws.Cells.Item(132, 14).Value = -22535

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 763.6
$ws.Cells.Item(100, 9).Value = 454.5
$ws.Cells.Item(100, 10).Value = 2000
$ws.Cells.Item(100, 11).Value = 909
$ws.Cells.Item(100, 12).Value = 4000
$ws.Cells.Item(100, 13).Value = -368

$ws.Cells.Item(107, 8).Value = 1392.4166
$ws.Cells.Item(107, 9).Value = 1425.4445
$ws.Cells.Item(107, 10).Value = 1293.3334
$ws.Cells.Item(107, 11).Value = 4276.333500000001
$ws.Cells.Item(107, 12).Value = 3880.0002
$ws.Cells.Item(107, 13).Value = -2356.333500000001
$ws.Cells.Item(107, 14).Value = -7720.0002

$ws.Cells.Item(109, 8).Value = 81659.664
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(109, 10).Value = 81659.664
$ws.Cells.Item(109, 11).Value = 0
$ws.Cells.Item(109, 12).Value = 81659.664
$ws.Cells.Item(109, 14).Value = -84433.664

$ws.Cells.Item(122, 8).Value = 74614.03
$ws.Cells.Item(122, 9).Value = 90529
$ws.Cells.Item(122, 10).Value = 5649.1665
$ws.Cells.Item(122, 11).Value = 271587
$ws.Cells.Item(122, 12).Value = 16947.4995
$ws.Cells.Item(122, 13).Value = -269137

$ws.Cells.Item(123, 8).Value = 80414
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 80414
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 80414
$ws.Cells.Item(123, 14).Value = -90214

$ws.Cells.Item(124, 8).Value = 77450
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 77450
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 77450
$ws.Cells.Item(124, 14).Value = -87270

$ws.Cells.Item(132, 8).Value = 16756.04
$ws.Cells.Item(132, 9).Value = 19363.28
$ws.Cells.Item(132, 10).Value = 2742.125
$ws.Cells.Item(132, 11).Value = 58089.84
$ws.Cells.Item(132, 12).Value = 8226.375
$ws.Cells.Item(132, 13).Value = -55559.84

$ws.Cells.Item(136, 8).Value = 23101.055
$ws.Cells.Item(136, 9).Value = 31259.309
$ws.Cells.Item(136, 10).Value = 3817.9092
$ws.Cells.Item(136, 11).Value = 93777.927
$ws.Cells.Item(136, 12).Value = 11453.7276
$ws.Cells.Item(136, 13).Value = -91227.927
